$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" ---
# This string is shared across the Overview sheet (E2, F2) and the per-locale
# status sheets (C2 on "zh-cn" and "de-de"). Update every occurrence so the
# underlying shared string is fully replaced (no cell is left pointing at the
# old text).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width changes ---
# Target stored OOXML width is 13.4101848602295 "character" units. Excel's
# ColumnWidth COM property is always quantized to whole-pixel increments
# (stored_width = round((ColumnWidth + 5/6) * 6) / 6), so we pick the
# ColumnWidth value whose quantized result lands closest to the target.
$newColumnWidth = 12.5

$wsOverview.Range("E1").ColumnWidth = $newColumnWidth
$wsOverview.Range("F1").ColumnWidth = $newColumnWidth

$wsZhCn.Range("C1").ColumnWidth = $newColumnWidth

$wsDeDe.Range("C1").ColumnWidth = $newColumnWidth
